# Apply the re-shuffled "Fecha / Volumen / Precio" figures for rows 2-14
# on the single worksheet of the workbook. Only columns D (Fecha), J
# (Volumen), K (Precio minimo), L (Precio maximo), M (Precio promedio
# ponderado) and P (Precio $/Kg) change; every other column keeps its
# current value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ D = fecha; J = volumen; K = precio min; L = precio max; M = precio prom; P = precio kg }
$rows = @{
    2  = @{ D = 44428; J = 300; K = 15000; L = 16000; M = 15500; P = 1550 }
    3  = @{ D = 44460; J = 300; K = 15000; L = 16000; M = 15500; P = 1550 }
    4  = @{ D = 44406; J = 400; K = 14000; L = 15000; M = 14500; P = 1450 }
    5  = @{ D = 44265; J = 200; K = 15000; L = 16000; M = 15500; P = 1550 }
    6  = @{ D = 44263; J = 300; K = 15000; L = 16000; M = 15500; P = 1550 }
    7  = @{ D = 44441; J = 300; K = 15000; L = 16000; M = 15500; P = 1550 }
    8  = @{ D = 44358; J = 300; K = 14000; L = 15000; M = 14500; P = 1450 }
    9  = @{ D = 44377; J = 650; K = 14000; L = 15000; M = 14538; P = 1454 }
    10 = @{ D = 44160; J = 360; K = 10000; L = 11000; M = 10500; P = 1050 }
    11 = @{ D = 44291; J = 200; K = 13000; L = 14000; M = 13500; P = 1350 }
    12 = @{ D = 44204; J = 400; K = 10000; L = 11000; M = 10500; P = 1050 }
    13 = @{ D = 44218; J = 320; K = 10000; L = 11000; M = 10500; P = 1050 }
    14 = @{ D = 44330; J = 300; K = 13000; L = 14000; M = 13500; P = 1350 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
